$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2-23) with new values per the diff ---
# Row 2
$ws.Range("D2").Value = 44488
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 25000
$ws.Range("O2").Value = 26000
$ws.Range("P2").Value = 25600
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("S2").Value = 2560
$ws.Range("T2").Value = 10

# Row 3
$ws.Range("D3").Value = 44511
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 25000
$ws.Range("O3").Value = 26000
$ws.Range("P3").Value = 25375
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 2538
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44160
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("Q4").Value = "$/bandeja 8 kilos"
$ws.Range("S4").Value = 2188
$ws.Range("T4").Value = 8

# Row 5
$ws.Range("D5").Value = 44160
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = "$/bandeja 8 kilos"
$ws.Range("S5").Value = 1875
$ws.Range("T5").Value = 8

# Row 7
$ws.Range("D7").Value = 44484
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 26000
$ws.Range("P7").Value = 25500
$ws.Range("Q7").Value = "$/bandeja 10 kilos"
$ws.Range("S7").Value = 2550
$ws.Range("T7").Value = 10

# Row 8
$ws.Range("D8").Value = 44516
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 1900
$ws.Range("O8").Value = 2000
$ws.Range("P8").Value = 1950
$ws.Range("S8").Value = 1950

# Row 9
$ws.Range("D9").Value = 44516
$ws.Range("L9").Value = "Segunda"
$ws.Range("N9").Value = 1700
$ws.Range("O9").Value = 1700
$ws.Range("P9").Value = 1700
$ws.Range("Q9").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S9").Value = 1700
$ws.Range("T9").Value = 1

# Row 10
$ws.Range("D10").Value = 44475
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 29000
$ws.Range("O10").Value = 30000
$ws.Range("P10").Value = 29500
$ws.Range("Q10").Value = "$/bandeja 10 kilos"
$ws.Range("S10").Value = 2950
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("L11").Value = "Primera"
$ws.Range("N11").Value = 2600
$ws.Range("O11").Value = 2600
$ws.Range("P11").Value = 2600
$ws.Range("S11").Value = 2600

# Row 12
$ws.Range("D12").Value = 44483
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 2400
$ws.Range("O12").Value = 2400
$ws.Range("P12").Value = 2400
$ws.Range("R12").Value = "Provincia de Limarí"
$ws.Range("S12").Value = 2400

# Row 13
$ws.Range("D13").Value = 44467
$ws.Range("N13").Value = 2700
$ws.Range("O13").Value = 2800
$ws.Range("P13").Value = 2750
$ws.Range("Q13").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S13").Value = 2750
$ws.Range("T13").Value = 1

# Row 14
$ws.Range("D14").Value = 44467
$ws.Range("N14").Value = 2500
$ws.Range("O14").Value = 2500
$ws.Range("P14").Value = 2500
$ws.Range("Q14").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S14").Value = 2500
$ws.Range("T14").Value = 1

# Row 15
$ws.Range("D15").Value = 44505
$ws.Range("N15").Value = 2200
$ws.Range("O15").Value = 2200
$ws.Range("P15").Value = 2200
$ws.Range("Q15").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S15").Value = 2200
$ws.Range("T15").Value = 1

# Row 16
$ws.Range("D16").Value = 44505
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 100
$ws.Range("N16").Value = 1800
$ws.Range("O16").Value = 1800
$ws.Range("P16").Value = 1800
$ws.Range("Q16").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S16").Value = 1800
$ws.Range("T16").Value = 1

# Row 17
$ws.Range("D17").Value = 44454
$ws.Range("N17").Value = 30000
$ws.Range("O17").Value = 31000
$ws.Range("P17").Value = 30500
$ws.Range("S17").Value = 3050

# Row 18
$ws.Range("D18").Value = 44469
$ws.Range("N18").Value = 28000
$ws.Range("O18").Value = 29000
$ws.Range("P18").Value = 28500
$ws.Range("Q18").Value = "$/bandeja 10 kilos"
$ws.Range("S18").Value = 2850
$ws.Range("T18").Value = 10

# Row 19
$ws.Range("D19").Value = 44446
$ws.Range("L19").Value = "Primera"
$ws.Range("N19").Value = 3200
$ws.Range("O19").Value = 3300
$ws.Range("P19").Value = 3250
$ws.Range("R19").Value = "Provincia del Elquí"
$ws.Range("S19").Value = 3250

# Row 20
$ws.Range("D20").Value = 44462
$ws.Range("N20").Value = 2900
$ws.Range("O20").Value = 3000
$ws.Range("P20").Value = 2950
$ws.Range("Q20").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S20").Value = 2950
$ws.Range("T20").Value = 1

# Row 21
$ws.Range("D21").Value = 44462
$ws.Range("L21").Value = "Segunda"
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 2600
$ws.Range("O21").Value = 2600
$ws.Range("P21").Value = 2600
$ws.Range("S21").Value = 2600

# Row 22
$ws.Range("D22").Value = 44491
$ws.Range("L22").Value = "Primera"
$ws.Range("M22").Value = 150
$ws.Range("N22").Value = 25000
$ws.Range("O22").Value = 26000
$ws.Range("P22").Value = 25467
$ws.Range("Q22").Value = "$/bandeja 10 kilos"
$ws.Range("S22").Value = 2547
$ws.Range("T22").Value = 10

# Row 23
$ws.Range("D23").Value = 44495
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 26000
$ws.Range("O23").Value = 27000
$ws.Range("P23").Value = 26500
$ws.Range("Q23").Value = "$/bandeja 10 kilos"
$ws.Range("R23").Value = "Provincia de Limarí"
$ws.Range("S23").Value = 2650
$ws.Range("T23").Value = 10

# --- Add new row 24 ---
$ws.Range("A24").Value = 11
$ws.Range("B24").Value = "Vega Monumental Concepción"
$ws.Range("C24").Value = "Bíobío"
$ws.Range("D24").Value = 44517
$ws.Range("E24").Value = 8
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100107
$ws.Range("H24").Value = "Otros"
$ws.Range("I24").Value = 100107002
$ws.Range("J24").Value = "Chirimoya"
$ws.Range("K24").Value = "Cultivar IV Región"
$ws.Range("L24").Value = "Primera"
$ws.Range("M24").Value = 100
$ws.Range("N24").Value = 25000
$ws.Range("O24").Value = 27000
$ws.Range("P24").Value = 26000
$ws.Range("Q24").Value = "$/bandeja 10 kilos"
$ws.Range("R24").Value = "Provincia de Limarí"
$ws.Range("S24").Value = 2600
$ws.Range("T24").Value = 10

# Match date style/number format used by other Fecha cells in column D
$ws.Range("D24").NumberFormat = $ws.Range("D23").NumberFormat
